$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column N (14th column) to host the new field
$ws.Columns.Item(14).Insert()

# Header for the new column
$ws.Cells.Item(1, 14).Value = "eIDAS RequesterID"

# Data value for the new column (numeric)
$ws.Cells.Item(2, 14).Value = 1234

# Update the view / selection to match the final state
$ws.Application.ActiveWindow.ScrollColumn = 10
$ws.Range("N3").Select()
